$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (PricePerHour), shifting it and every
# column after it one slot to the right.
$ws.Range("C:C").Insert()

# New header: "AvailableSpots" goes between "TotalSpots" (B) and the shifted
# "PricePerHour" (now D).
$ws.Range("C1").Value = "AvailableSpots"

# AvailableSpots starts out equal to TotalSpots for every data row, and the
# (now-shifted) PricePerHour column resets to 0. Cells in this sheet carry a
# Text ("@") number format, so a plain `.Value = <number>` assignment would
# get stored as text; temporarily clearing the style keeps the write numeric,
# then the original Text format is restored to match the sheet's styling.
for ($r = 2; $r -le 6; $r++) {
    $totalSpots = $ws.Cells.Item($r, 2).Formula

    $available = $ws.Cells.Item($r, 3)
    $available.Style = "Normal"
    $available.Value = $totalSpots
    $available.NumberFormat = "@"

    $price = $ws.Cells.Item($r, 4)
    $price.Style = "Normal"
    $price.Value = 0
    $price.NumberFormat = "@"
}

# Move the active selection, matching the post-edit cursor position.
$ws.Range("A7").Select()
